$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.34889829158783
$ws.Range("B1").Value = 1.532660007476807
$ws.Range("C1").Value = 1.887895941734314
$ws.Range("D1").Value = 2.557666778564453
$ws.Range("E1").Value = 6.629444122314453
